{"js": "// UC-06 \"Busquedas Libro\" textual use-case update.\n//\n// 1) Drop the obsolete \"1- El administrador abre el m\u00f3dulo de Search\" step.\n// 2) Renumber/reword the remaining \"b\u00e1sico\" steps (2.x -> 1.x) and justify them.\n// 3) Justify the \"Escenario Alternativo(B\u00fasquedas)\" heading.\n// 4) Renumber/reword the \"alternativo\" steps (2.2a/1-3 -> 1.2 a/1-3) and justify them.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate + delete the \"Search\" module paragraph outright (old step 1 of the\n// basic scenario gets removed entirely when the steps were renumbered).\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"El administrador abre el m\u00f3dulo de Search\") !== -1) {\n    items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n\n// Re-load paragraphs now that the document shifted.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\nconst items2 = paragraphs2.items;\n\n// Old text (exact match) -> new text, for the basic-scenario cell.\nconst basicRewrites = [\n  [\n    \"2- El administrador selecciona la  para el tipo de b\u00fasqueda que desea realizar(T\u00edtulos, Autores, Encabezado Materia, Clasificaciones, Series, Ejemplares, Usuarios, Localizaciones )\",\n    \"1- El administrador selecciona la  para el tipo de b\u00fasqueda que desea realizar ya sea por t\u00edtulos, autores, encabezado materia, clasificaciones, series, ejemplares, usuarios, localizaciones.\"\n  ],\n  [\n    \"    2.1- El administrador escribe en textbox  el texto a buscar\",\n    \"    1.1- El administrador digita lo que desea buscar.\"\n  ],\n  [\n    \"    2.2- El sistema consulta la base de datos buscando coincidencias\",\n    \"    1.2- El sistema consulta buscando coincidencias.\"\n  ],\n  [\n    \"    2.3- El sistema muestra en un datagrid la lista de los elementos que contengan coincidencia       con el texto ingresado\",\n    \"    1.3- El sistema muestra una lista de los elementos que contengan coincidencia con el texto ingresado.\"\n  ],\n  [\n    \"    2.4- El administrador le da doble click sobre el elemento\",\n    \"    1.4- El administrador selecciona el elemento.\"\n  ],\n  [\n    \"    2.5- El sistema muestra toda la informaci\u00f3n referente al elemento seleccionado\",\n    \"    1.5- El sistema muestra toda la informaci\u00f3n referente al elemento seleccionado.\"\n  ]\n];\n\n// Old text (exact match) -> new text, for the heading + alternativo cell.\nconst altRewrites = [\n  [\n    \"Escenario Alternativo(B\u00fasquedas)\",\n    null // text unchanged, justification only\n  ],\n  [\n    \"2.2a- El sistema no encuentra coincidencias \",\n    \"1.2 a- El sistema no encuentra coincidencias.\"\n  ],\n  [\n    \"     1- El sistema envia un mensaje de informaci\u00f3n indicando que el libro no existe\",\n    \"     1- El sistema env\u00eda un mensaje de informaci\u00f3n indicando que el libro no existe.\"\n  ],\n  [\n    \"     2- El sistema muestra un datagrid vac\u00edo\",\n    \"     2- El sistema muestra una lista vac\u00eda.\"\n  ],\n  [\n    \"     3- El sistema no conecta con la base de datos y muestra un mensaje de error de conexi\u00f3n\",\n    \"     3- El sistema no se conecta y muestra un mensaje de error de conexi\u00f3n.\"\n  ],\n  [\n    \"    4- El administrador vuelve al paso 2.1\",\n    null // text unchanged, justification only\n  ]\n];\n\nconst allRewrites = basicRewrites.concat(altRewrites);\n\nfor (let i = 0; i < items2.length; i++) {\n  const para = items2[i];\n  const text = para.text;\n  for (let j = 0; j < allRewrites.length; j++) {\n    const oldText = allRewrites[j][0];\n    const newText = allRewrites[j][1];\n    if (text === oldText) {\n      if (newText !== null) {\n        para.insertText(newText, Word.InsertLocation.replace);\n      }\n      para.alignment = Word.Alignment.justified;\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# UC-06 \"Busquedas Libro\" textual use-case update.\n#\n# 1) Drop the obsolete \"1- El administrador abre el m\u00f3dulo de Search\" step.\n# 2) Renumber/reword the remaining \"b\u00e1sico\" steps (2.x -> 1.x) and justify them.\n# 3) Justify the \"Escenario Alternativo(B\u00fasquedas)\" heading.\n# 4) Renumber/reword the \"alternativo\" steps (2.2a/1-3 -> 1.2 a/1-3) and justify them.\n\n$d = $word.ActiveDocument\n$table = $d.Tables(1)\n\n$wdAlignParagraphJustify = 3\n\n# Paragraph.Range / the Find-and-Replace APIs only touch the first run of a\n# paragraph when later runs are separated by non-run siblings (e.g.\n# <w:proofErr/> around a flagged word), so replace paragraph text by building\n# an explicit range that spans the whole paragraph *except* its trailing\n# paragraph mark (Range.End - 1) before assigning .Text. That guarantees all\n# runs in the paragraph are replaced, not just the first one.\nfunction Set-ParagraphText($paragraph, $newText) {\n    $r = $paragraph.Range\n    $full = $d.Range($r.Start, $r.End - 1)\n    $full.Text = $newText\n}\n\n# --- Row 2: basic-scenario steps -------------------------------------------\n$basicCell = $table.Cell(2, 1)\n\n# Old step 1 (\"... abre el m\u00f3dulo de Search\") is removed outright; the other\n# steps shift up and get renumbered from 2.x to 1.x.\n$basicCell.Range.Paragraphs(1).Range.Delete()\n\nSet-ParagraphText $basicCell.Range.Paragraphs(1) \"1- El administrador selecciona la  para el tipo de b\u00fasqueda que desea realizar ya sea por t\u00edtulos, autores, encabezado materia, clasificaciones, series, ejemplares, usuarios, localizaciones.\"\nSet-ParagraphText $basicCell.Range.Paragraphs(2) \"    1.1- El administrador digita lo que desea buscar.\"\nSet-ParagraphText $basicCell.Range.Paragraphs(3) \"    1.2- El sistema consulta buscando coincidencias.\"\nSet-ParagraphText $basicCell.Range.Paragraphs(4) \"    1.3- El sistema muestra una lista de los elementos que contengan coincidencia con el texto ingresado.\"\nSet-ParagraphText $basicCell.Range.Paragraphs(5) \"    1.4- El administrador selecciona el elemento.\"\nSet-ParagraphText $basicCell.Range.Paragraphs(6) \"    1.5- El sistema muestra toda la informaci\u00f3n referente al elemento seleccionado.\"\n\nfor ($i = 1; $i -le 6; $i++) {\n    $basicCell.Range.Paragraphs($i).Alignment = $wdAlignParagraphJustify\n}\n\n# --- Row 3: \"Escenario Alternativo(B\u00fasquedas)\" heading ---------------------\n$headingCell = $table.Cell(3, 1)\n$headingCell.Range.Paragraphs(1).Alignment = $wdAlignParagraphJustify\n\n# --- Row 4: alternative-scenario steps --------------------------------------\n$altCell = $table.Cell(4, 1)\n\nSet-ParagraphText $altCell.Range.Paragraphs(1) \"1.2 a- El sistema no encuentra coincidencias.\"\nSet-ParagraphText $altCell.Range.Paragraphs(2) \"     1- El sistema env\u00eda un mensaje de informaci\u00f3n indicando que el libro no existe.\"\nSet-ParagraphText $altCell.Range.Paragraphs(3) \"     2- El sistema muestra una lista vac\u00eda.\"\nSet-ParagraphText $altCell.Range.Paragraphs(4) \"     3- El sistema no se conecta y muestra un mensaje de error de conexi\u00f3n.\"\n# Paragraph 5 (\"4- El administrador vuelve al paso 2.1\") keeps its text.\n\nfor ($i = 1; $i -le 5; $i++) {\n    $altCell.Range.Paragraphs($i).Alignment = $wdAlignParagraphJustify\n}\n"}
